$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E13").Value = "SMMartin Merisalu"
